$wb = $excel.ActiveWorkbook

# --- RB sheet: add Week 15 simulation row for C.Reynolds ---
$wsRB = $wb.Worksheets.Item("RB")
$wsRB.Cells.Item(7, 1).Value = "C.Reynolds"
for ($c = 2; $c -le 10; $c++) {
    $wsRB.Cells.Item(7, $c).Value = 0
}
$wsRB.Range("I8").Select()

# --- TE sheet: add Week 15 simulation row for S.Zylstra ---
$wsTE = $wb.Worksheets.Item("TE")
$wsTE.Cells.Item(6, 1).Value = "S.Zylstra"
for ($c = 2; $c -le 10; $c++) {
    $wsTE.Cells.Item(6, $c).Value = 0
}
$wsTE.Range("J7").Select()

# TE becomes the active/selected tab
$wsTE.Activate()
